$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new match row: Grigor Dimitrov def. Roger Federer, Australian Open Final (2018-01-23)
$ws.Range("A9").Value = 43123
$ws.Range("A9").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B9").Value = "Grigor Dimitrov"
$ws.Range("C9").Value = "Roger Federer"
$ws.Range("D9").Value = "6–4, 3–6, 6–1, 3–6, 6–3"
$ws.Range("E9").Value = "Australian Open"
$ws.Range("F9").Value = "F"

$ws.Columns.Item(1).ColumnWidth = 12.7

$ws.Range("F9").Select() | Out-Null
